$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.20127533333333
$ws.Range("H2").Value = 33.603826
$ws.Range("I2").Value = 0.1186573945858706
$ws.Range("J2").Value = 0.1186573945858706
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 9.460162333333333
$ws.Range("N2").Value = 28.380487
$ws.Range("O2").Value = 0.08197024919772995
$ws.Range("P2").Value = 0.08197024919772995
$ws.Range("Q2").Value = 105.9658829936958
$ws.Range("R2").Value = 953.6929469432619
$ws.Range("S2").Value = 0.009726376203357189
$ws.Range("T2").Value = 0.009726376203357189

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.20127533333333
$ws.Range("H3").Value = 33.603826
$ws.Range("I3").Value = 0.1186573945858706
$ws.Range("J3").Value = 0.1186573945858706
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.18256633333333
$ws.Range("N3").Value = 138.547699
$ws.Range("O3").Value = 0.4001618933742075
$ws.Range("P3").Value = 0.4001618933742075
$ws.Range("Q3").Value = 517.3036410995971
$ws.Range("R3").Value = 4655.732769896374
$ws.Range("S3").Value = 0.04748216768033243
$ws.Range("T3").Value = 0.04748216768033243

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.20127533333333
$ws.Range("H4").Value = 33.603826
$ws.Range("I4").Value = 0.1186573945858706
$ws.Range("J4").Value = 0.1186573945858706
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 51.40166966666666
$ws.Range("N4").Value = 154.205009
$ws.Range("O4").Value = 0.4453842886934318
$ws.Range("P4").Value = 0.4453842886934319
$ws.Range("Q4").Value = 575.7642545293816
$ws.Range("R4").Value = 5181.878290764434
$ws.Range("S4").Value = 0.05284813928584386
$ws.Range("T4").Value = 0.05284813928584386

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 11.20127533333333
$ws.Range("H5").Value = 33.603826
$ws.Range("I5").Value = 0.1186573945858706
$ws.Range("J5").Value = 0.1186573945858706
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.365307333333334
$ws.Range("N5").Value = 25.095922
$ws.Range("O5").Value = 0.07248356873463072
$ws.Range("P5").Value = 0.07248356873463073
$ws.Range("Q5").Value = 93.70211068861911
$ws.Range("R5").Value = 843.318996197572
$ws.Range("S5").Value = 0.008600711416337153
$ws.Range("T5").Value = 0.008600711416337154

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 47.94465366666667
$ws.Range("H6").Value = 143.833961
$ws.Range("I6").Value = 0.5078874966566524
$ws.Range("J6").Value = 0.5078874966566524
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 9.460162333333333
$ws.Range("N6").Value = 28.380487
$ws.Range("O6").Value = 0.08197024919772995
$ws.Range("P6").Value = 0.08197024919772995
$ws.Range("Q6").Value = 453.564206702112
$ws.Range("R6").Value = 4082.077860319007
$ws.Range("S6").Value = 0.04163166466535704
$ws.Range("T6").Value = 0.04163166466535704

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 47.94465366666667
$ws.Range("H7").Value = 143.833961
$ws.Range("I7").Value = 0.5078874966566524
$ws.Range("J7").Value = 0.5078874966566524
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.18256633333333
$ws.Range("N7").Value = 138.547699
$ws.Range("O7").Value = 0.4001618933742075
$ws.Range("P7").Value = 0.4001618933742075
$ws.Range("Q7").Value = 2214.207148289527
$ws.Range("R7").Value = 19927.86433460574
$ws.Range("S7").Value = 0.2032372222832125
$ws.Range("T7").Value = 0.2032372222832125

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 47.94465366666667
$ws.Range("H8").Value = 143.833961
$ws.Range("I8").Value = 0.5078874966566524
$ws.Range("J8").Value = 0.5078874966566524
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 51.40166966666666
$ws.Range("N8").Value = 154.205009
$ws.Range("O8").Value = 0.4453842886934318
$ws.Range("P8").Value = 0.4453842886934319
$ws.Range("Q8").Value = 2464.435250056739
$ws.Range("R8").Value = 22179.91725051065
$ws.Range("S8").Value = 0.2262051114347109
$ws.Range("T8").Value = 0.2262051114347109

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 47.94465366666667
$ws.Range("H9").Value = 143.833961
$ws.Range("I9").Value = 0.5078874966566524
$ws.Range("J9").Value = 0.5078874966566524
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.365307333333334
$ws.Range("N9").Value = 25.095922
$ws.Range("O9").Value = 0.07248356873463072
$ws.Range("P9").Value = 0.07248356873463073
$ws.Range("Q9").Value = 401.0717629118936
$ws.Range("R9").Value = 3609.645866207043
$ws.Range("S9").Value = 0.036813498273372
$ws.Range("T9").Value = 0.03681349827337201

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 25.59984766666667
$ws.Range("H10").Value = 76.799543
$ws.Range("I10").Value = 0.2711844085184091
$ws.Range("J10").Value = 0.2711844085184091
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 9.460162333333333
$ws.Range("N10").Value = 28.380487
$ws.Range("O10").Value = 0.08197024919772995
$ws.Range("P10").Value = 0.08197024919772995
$ws.Range("Q10").Value = 242.1787146352712
$ws.Range("R10").Value = 2179.608431717441
$ws.Range("S10").Value = 0.022229053544793
$ws.Range("T10").Value = 0.022229053544793

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.59984766666667
$ws.Range("H11").Value = 76.799543
$ws.Range("I11").Value = 0.2711844085184091
$ws.Range("J11").Value = 0.2711844085184091
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.18256633333333
$ws.Range("N11").Value = 138.547699
$ws.Range("O11").Value = 0.4001618933742075
$ws.Range("P11").Value = 0.4001618933742075
$ws.Range("Q11").Value = 1182.266662989062
$ws.Range("R11").Value = 10640.39996690156
$ws.Range("S11").Value = 0.1085176663662912
$ws.Range("T11").Value = 0.1085176663662912

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.59984766666667
$ws.Range("H12").Value = 76.799543
$ws.Range("I12").Value = 0.2711844085184091
$ws.Range("J12").Value = 0.2711844085184091
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 51.40166966666666
$ws.Range("N12").Value = 154.205009
$ws.Range("O12").Value = 0.4453842886934318
$ws.Range("P12").Value = 0.4453842886934319
$ws.Range("Q12").Value = 1315.874913278987
$ws.Range("R12").Value = 11842.87421951089
$ws.Range("S12").Value = 0.1207812748927207
$ws.Range("T12").Value = 0.1207812748927207

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.59984766666667
$ws.Range("H13").Value = 76.799543
$ws.Range("I13").Value = 0.2711844085184091
$ws.Range("J13").Value = 0.2711844085184091
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 8.365307333333334
$ws.Range("N13").Value = 25.095922
$ws.Range("O13").Value = 0.07248356873463072
$ws.Range("P13").Value = 0.07248356873463073
$ws.Range("Q13").Value = 214.1505934181829
$ws.Range("R13").Value = 1927.355340763646
$ws.Range("S13").Value = 0.01965641371460428
$ws.Range("T13").Value = 0.01965641371460429

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 9.654369000000001
$ws.Range("H14").Value = 28.963107
$ws.Range("I14").Value = 0.1022707002390678
$ws.Range("J14").Value = 0.1022707002390678
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 9.460162333333333
$ws.Range("N14").Value = 28.380487
$ws.Range("O14").Value = 0.08197024919772995
$ws.Range("P14").Value = 0.08197024919772995
$ws.Range("Q14").Value = 91.33189796590101
$ws.Range("R14").Value = 821.987081693109
$ws.Range("S14").Value = 0.008383154784222727
$ws.Range("T14").Value = 0.008383154784222725

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 9.654369000000001
$ws.Range("H15").Value = 28.963107
$ws.Range("I15").Value = 0.1022707002390678
$ws.Range("J15").Value = 0.1022707002390678
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 46.18256633333333
$ws.Range("N15").Value = 138.547699
$ws.Range("O15").Value = 0.4001618933742075
$ws.Range("P15").Value = 0.4001618933742075
$ws.Range("Q15").Value = 445.863536748977
$ws.Range("R15").Value = 4012.771830740793
$ws.Range("S15").Value = 0.04092483704437138
$ws.Range("T15").Value = 0.04092483704437137

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 9.654369000000001
$ws.Range("H16").Value = 28.963107
$ws.Range("I16").Value = 0.1022707002390678
$ws.Range("J16").Value = 0.1022707002390678
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 51.40166966666666
$ws.Range("N16").Value = 154.205009
$ws.Range("O16").Value = 0.4453842886934318
$ws.Range("P16").Value = 0.4453842886934319
$ws.Range("Q16").Value = 496.250686178107
$ws.Range("R16").Value = 4466.256175602963
$ws.Range("S16").Value = 0.0455497630801564
$ws.Range("T16").Value = 0.0455497630801564

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 9.654369000000001
$ws.Range("H17").Value = 28.963107
$ws.Range("I17").Value = 0.1022707002390678
$ws.Range("J17").Value = 0.1022707002390678
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 8.365307333333334
$ws.Range("N17").Value = 25.095922
$ws.Range("O17").Value = 0.07248356873463072
$ws.Range("P17").Value = 0.07248356873463073
$ws.Range("Q17").Value = 80.76176379440601
$ws.Range("R17").Value = 726.855874149654
$ws.Range("S17").Value = 0.007412945330317284
$ws.Range("T17").Value = 0.007412945330317285
